$d = $word.ActiveDocument

# --- Edit 1: "Date: dd-mm-yyyy" paragraph -> split into two runs around
#     a <w:proofErr spellStart/spellEnd> pair, and drop the _GoBack bookmark.
$dateTarget = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Date: dd-mm-yyyy") {
        $dateTarget = $p
        break
    }
}
if ($dateTarget -eq $null) { throw "Could not locate the 'Date: dd-mm-yyyy' paragraph" }
$dateRange = $dateTarget.Range
$dateXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Date: dd-mm-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>yyyy</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$dateRange.InsertXML($dateXml)

# --- Edit 2: "[file:8435.jpg]" paragraph -> split into three runs around
#     a <w:proofErr gramStart/gramEnd> pair wrapping "file:8435.jpg".
$fileTarget = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "[file:8435.jpg]") {
        $fileTarget = $p
        break
    }
}
if ($fileTarget -eq $null) { throw "Could not locate the '[file:8435.jpg]' paragraph" }
$fileRange = $fileTarget.Range
$fileXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>[</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>file:8435.jpg</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>]</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$fileRange.InsertXML($fileXml)

# --- Edit 3: append a brand-new QN=2 question table after the existing
#     table, followed by a fresh trailing empty paragraph before sectPr.
$lastTbl = $d.Tables(1)
$anchorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $lastTbl.Range.End) {
        $anchorPara = $p
        break
    }
}
if ($anchorPara -eq $null) { throw "Could not locate the paragraph right after the existing table" }
$tableXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="10000" w:type="dxa"/><w:tblLook w:val="01E0" w:firstRow="1" w:lastRow="1" w:firstColumn="1" w:lastColumn="1" w:noHBand="0" w:noVBand="0"/></w:tblPr><w:tblGrid><w:gridCol w:w="2000"/><w:gridCol w:w="8000"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="2000" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:t>Q</w:t></w:r><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>N=2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="8000" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>10</w:t></w:r><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>+</w:t></w:r><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>10</w:t></w:r><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve"> = ?</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>a.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="8000" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>1</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>b.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="8000" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>2</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>c.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="8000" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>3</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>d.</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="8000" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>0</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>ANSWER:</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="8000" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>D</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>MARK:</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="8000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>0.5</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>UNIT:</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="8000" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:t>Chapter</w:t></w:r><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>2</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2000" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>MIX CHOICES:</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="8000" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:rPr><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="vi-VN"/></w:rPr><w:t>No</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p><w:r><w:t>ZZZPLACEHOLDERZZZ</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$anchorPara.Range.InsertXML($tableXml)

# The inserted trailing paragraph carries placeholder text so the engine
# does not collapse/merge it away as a redundant empty paragraph right
# before </w:body>'s <w:sectPr>; now blank out the placeholder text while
# leaving the (now ordinary) empty paragraph mark intact.
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalRange = $finalPara.Range
$finalRange.MoveEnd(1, -1)
$finalRange.Delete()
